$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 4804
$ws.Range("E2").Value = 279
$ws.Range("F2").Value = 279
$ws.Range("G2").Value = 278
$ws.Range("H2").Value = 163
$ws.Range("I2").Value = 163
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 2072
$ws.Range("L2").Value = 872
$ws.Range("M2").Value = 1200
$ws.Range("N2").Value = 1192
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 25
$ws.Range("Q2").Value = 481
$ws.Range("R2").Value = -421
$ws.Range("S2").Value = -267
$ws.Range("T2").Value = 398
$ws.Range("U2").Value = 84
$ws.Range("V2").Value = 410
$ws.Range("W2").Value = 5.81
$ws.Range("X2").Value = 3.4
$ws.Range("Y2").Value = 13.3
$ws.Range("Z2").Value = 7.26
$ws.Range("AA2").Value = 72.7
$ws.Range("AB2").Value = 4603.42
$ws.Range("AC2").Value = 326
$ws.Range("AD2").ClearContents()
$ws.Range("AE2").Value = 2383
$ws.Range("AF2").Value = 0
$ws.Range("AH2").ClearContents()
$ws.Range("AI2").Value = 122.75
$ws.Range("AJ2").Value = 50041193

# Row 3
$ws.Range("D3").Value = 3721
$ws.Range("E3").Value = -135
$ws.Range("F3").Value = -135
$ws.Range("G3").Value = -164
$ws.Range("H3").Value = -148
$ws.Range("I3").Value = -145
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 2471
$ws.Range("L3").Value = 1600
$ws.Range("M3").Value = 870
$ws.Range("N3").Value = 860
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 25
$ws.Range("Q3").Value = -71
$ws.Range("R3").Value = -314
$ws.Range("S3").Value = 349
$ws.Range("T3").Value = 464
$ws.Range("U3").Value = -535
$ws.Range("V3").Value = 953
$ws.Range("W3").Value = -3.62
$ws.Range("X3").Value = -3.98
$ws.Range("Y3").Value = -14.09
$ws.Range("Z3").Value = -6.51
$ws.Range("AA3").Value = 183.88
$ws.Range("AB3").Value = 3247.5
$ws.Range("AC3").Value = -289
$ws.Range("AD3").ClearContents()
$ws.Range("AE3").Value = 1719
$ws.Range("AF3").Value = 0
$ws.Range("AG3").Value = 0
$ws.Range("AH3").ClearContents()
$ws.Range("AI3").Value = 0
$ws.Range("AJ3").Value = 50041193

# Row 4
$ws.Range("D4").Value = 5572
$ws.Range("E4").Value = 42
$ws.Range("F4").Value = 42
$ws.Range("G4").Value = -1
$ws.Range("H4").Value = -4
$ws.Range("I4").Value = -3
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 2684
$ws.Range("L4").Value = 1823
$ws.Range("M4").Value = 861
$ws.Range("N4").Value = 853
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 25
$ws.Range("Q4").Value = 306
$ws.Range("R4").Value = -157
$ws.Range("S4").Value = -79
$ws.Range("T4").Value = 219
$ws.Range("U4").Value = 87
$ws.Range("V4").Value = 877
$ws.Range("W4").Value = 0.76
$ws.Range("X4").Value = -0.08
$ws.Range("Y4").Value = -0.33
$ws.Range("Z4").Value = -0.17
$ws.Range("AA4").Value = 211.59
$ws.Range("AB4").Value = 3254.56
$ws.Range("AC4").Value = -6
$ws.Range("AD4").ClearContents()
$ws.Range("AE4").Value = 1705
$ws.Range("AF4").Value = 0
$ws.Range("AG4").Value = 0
$ws.Range("AH4").ClearContents()
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = 50041193

# Row 5
$ws.Range("D5").Value = 6794
$ws.Range("E5").Value = 406
$ws.Range("F5").Value = 406
$ws.Range("G5").Value = 314
$ws.Range("H5").Value = 249
$ws.Range("I5").Value = 251
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 2805
$ws.Range("L5").Value = 1710
$ws.Range("M5").Value = 1095
$ws.Range("N5").Value = 1089
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 25
$ws.Range("Q5").Value = 236
$ws.Range("R5").Value = -280
$ws.Range("S5").Value = 11
$ws.Range("T5").Value = 309
$ws.Range("U5").Value = -72
$ws.Range("V5").Value = 835
$ws.Range("W5").Value = 5.97
$ws.Range("X5").Value = 3.67
$ws.Range("Y5").Value = 25.88
$ws.Range("Z5").Value = 9.08
$ws.Range("AA5").Value = 156.08
$ws.Range("AB5").Value = 4232.08
$ws.Range("AC5").Value = 501
$ws.Range("AD5").ClearContents()
$ws.Range("AE5").Value = 2167
$ws.Range("AF5").Value = 0
$ws.Range("AG5").Value = 0
$ws.Range("AH5").ClearContents()
$ws.Range("AI5").Value = 0
$ws.Range("AJ5").Value = 50252570

# Row 6
$ws.Range("D6").Value = 6015
$ws.Range("E6").Value = 424
$ws.Range("F6").Value = 424
$ws.Range("G6").Value = 406
$ws.Range("H6").Value = 371
$ws.Range("I6").Value = 370
$ws.Range("K6").Value = 2849
$ws.Range("L6").Value = 1422
$ws.Range("M6").Value = 1428
$ws.Range("N6").Value = 1421
$ws.Range("P6").Value = 25
$ws.Range("Q6").Value = 621
$ws.Range("R6").Value = -476
$ws.Range("S6").Value = -54
$ws.Range("T6").Value = 416
$ws.Range("U6").Value = 206
$ws.Range("V6").Value = 801
$ws.Range("W6").Value = 7.05
$ws.Range("X6").Value = 6.16
$ws.Range("Y6").Value = 29.48
$ws.Range("Z6").Value = 13.11
$ws.Range("AA6").Value = 99.58
$ws.Range("AB6").Value = 5682.12
$ws.Range("AC6").Value = 736
$ws.Range("AD6").ClearContents()
$ws.Range("AE6").Value = 2827
$ws.Range("AF6").Value = 0
$ws.Range("AG6").Value = 0
$ws.Range("AH6").ClearContents()
$ws.Range("AI6").Value = 0
$ws.Range("AJ6").Value = 50252570

# Row 7
$ws.Range("D7").Value = 5961
$ws.Range("E7").Value = 331
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").Value = 279
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").Value = 5.55
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").Value = 488
$ws.Range("AD7").Value = 12.49
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()

# Row 8
$ws.Range("D8").Value = 12425
$ws.Range("E8").Value = 741
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").Value = 395
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").Value = 5.96
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").Value = 647
$ws.Range("AD8").Value = 9.41
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()

# Row 9
$ws.Range("D9").Value = 14451
$ws.Range("E9").Value = 913
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").Value = 490
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").Value = 6.32
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").Value = 802
$ws.Range("AD9").Value = 7.59
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()
